# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value2 = "61.244.23"
$cD.Style = "Normal"
$ws.Range("E2").Value2 = "  +0.74%  "

$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value2 = "2.928.74"
$cD.Style = "Normal"
$ws.Range("E3").Value2 = "  +0.55%  "

$cD = $ws.Range("D4")
$cD.NumberFormat = "@"
$cD.Value2 = "1.00"
$cD.Style = "Normal"
$ws.Range("E4").Value2 = "  -0.06%  "

$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value2 = "590.84"
$cD.Style = "Normal"
$ws.Range("E5").Value2 = "  +1.05%  "

$cD = $ws.Range("D6")
$cD.NumberFormat = "@"
$cD.Value2 = "146.21"
$cD.Style = "Normal"
$ws.Range("E6").Value2 = "  -0.32%  "

$cD = $ws.Range("D7")
$cD.NumberFormat = "@"
$cD.Value2 = "1.00"
$cD.Style = "Normal"
$ws.Range("E7").Value2 = "  -0.10%  "

$ws.Range("E8").Value2 = "  +1.50%  "

$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value2 = "2.927.73"
$cD.Style = "Normal"
$ws.Range("E9").Value2 = "  +0.56%  "

$cD = $ws.Range("D10")
$cD.NumberFormat = "@"
$cD.Value2 = "6.82"
$cD.Style = "Normal"
$ws.Range("E10").Value2 = "  +1.58%  "

$ws.Range("E11").Value2 = "  +1.17%  "

$ws.Range("E12").Value2 = "  -0.54%  "

$ws.Range("E13").Value2 = "  +2.14%  "

$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value2 = "33.85"
$cD.Style = "Normal"
$ws.Range("E14").Value2 = "  -1.34%  "

$ws.Range("E15").Value2 = "  +0.08%  "

$cD = $ws.Range("D16")
$cD.NumberFormat = "@"
$cD.Value2 = "3.413.35"
$cD.Style = "Normal"
$ws.Range("E16").Value2 = "  +0.35%  "

$cD = $ws.Range("D17")
$cD.NumberFormat = "@"
$cD.Value2 = "61.210.60"
$cD.Style = "Normal"
$ws.Range("E17").Value2 = "  +0.59%  "

$ws.Range("E18").Value2 = "  -0.96%  "

$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value2 = "2.926.76"
$cD.Style = "Normal"
$ws.Range("E19").Value2 = "  +0.14%  "

$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value2 = "432.21"
$cD.Style = "Normal"
$ws.Range("E20").Value2 = "  +1.77%  "

$ws.Range("E21").Value2 = "  -0.82%  "

$ws.Range("E22").Value2 = "  +2.36%  "

$ws.Range("E23").Value2 = "  -0.69%  "

$ws.Range("E24").Value2 = "  +0.21%  "

$cD = $ws.Range("D25")
$cD.NumberFormat = "@"
$cD.Value2 = "11.06"
$cD.Style = "Normal"
$ws.Range("E25").Value2 = "  +0.59%  "

$ws.Range("E26").Value2 = "  +3.79%  "

$cD = $ws.Range("D27")
$cD.NumberFormat = "@"
$cD.Value2 = "12.09"
$cD.Style = "Normal"
$ws.Range("E27").Value2 = "  +2.29%  "

$ws.Range("E28").Value2 = "  +0.17%  "

$cD = $ws.Range("D29")
$cD.NumberFormat = "@"
$cD.Value2 = "2.33"
$cD.Style = "Normal"
$ws.Range("E29").Value2 = "  +7.68%  "

$cD = $ws.Range("D30")
$cD.NumberFormat = "@"
$cD.Value2 = "1.00"
$cD.Style = "Normal"
$ws.Range("E30").Value2 = "  -0.10%  "

$ws.Range("E31").Value2 = "  +0.35%  "

$cD = $ws.Range("D32")
$cD.NumberFormat = "@"
$cD.Value2 = "7.18"
$cD.Style = "Normal"
$ws.Range("E32").Value2 = "  -0.60%  "

$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value2 = "26.59"
$cD.Style = "Normal"
$ws.Range("E33").Value2 = "  -0.12%  "

$ws.Range("E34").Value2 = "  +2.44%  "

$cD = $ws.Range("D35")
$cD.NumberFormat = "@"
$cD.Value2 = "0.0₃0871"
$cD.Style = "Normal"
$ws.Range("E35").Value2 = "  +3.78%  "

$ws.Range("E36").Value2 = "  +0.53%  "

$cD = $ws.Range("D37")
$cD.NumberFormat = "@"
$cD.Value2 = "3.12"
$cD.Style = "Normal"
$ws.Range("E37").Value2 = "  +3.36%  "

$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value2 = "5.66"
$cD.Style = "Normal"
$ws.Range("E38").Value2 = "  +0.08%  "

$cD = $ws.Range("D39")
$cD.NumberFormat = "@"
$cD.Value2 = "49.97"
$cD.Style = "Normal"
$ws.Range("E39").Value2 = "  +0.36%  "

$ws.Range("E40").Value2 = "  -0.27%  "

$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value2 = "0.124"
$cD.Style = "Normal"
$ws.Range("E41").Value2 = "  +0.93%  "

$cD = $ws.Range("D42")
$cD.NumberFormat = "@"
$cD.Value2 = "8.61"
$cD.Style = "Normal"
$ws.Range("E42").Value2 = "  -1.67%  "

$cD = $ws.Range("D43")
$cD.NumberFormat = "@"
$cD.Value2 = "0.290"
$cD.Style = "Normal"
$ws.Range("E43").Value2 = "  +1.12%  "

$cD = $ws.Range("D44")
$cD.NumberFormat = "@"
$cD.Value2 = "39.94"
$cD.Style = "Normal"
$ws.Range("E44").Value2 = "  -4.62%  "

$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value2 = "381.14"
$cD.Style = "Normal"
$ws.Range("E45").Value2 = "  +1.94%  "

$ws.Range("E46").Value2 = "  +0.94%  "

$cD = $ws.Range("D47")
$cD.NumberFormat = "@"
$cD.Value2 = "2.715.84"
$cD.Style = "Normal"
$ws.Range("E47").Value2 = "  +2.31%  "

$cD = $ws.Range("D48")
$cD.NumberFormat = "@"
$cD.Value2 = "130.90"
$cD.Style = "Normal"
$ws.Range("E48").Value2 = "  -1.31%  "

$cD = $ws.Range("D50")
$cD.NumberFormat = "@"
$cD.Value2 = "24.30"
$cD.Style = "Normal"
$ws.Range("E50").Value2 = "  -4.27%  "

$ws.Range("E51").Value2 = "  +0.58%  "

